# Apply the price/volume refresh for cryptos.xlsx as captured by the
# "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.398.08'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.48%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.655.94'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.60%  '

$ws.Range("E4").Value = '  +0.26%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.13%  '

$ws.Range("E6").Value = '  +0.17%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3917'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.49%  '

$ws.Range("E8").Value = '  -2.69%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.002'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.38%  '

$ws.Range("E10").Value = '  -6.18%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '50.15'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.04%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08534'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.09%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.90'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.24%  '

$ws.Range("E14").Value = '  -3.64%  '

$ws.Range("E15").Value = '  -3.06%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.623'
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.658.62'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.19%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '93.61'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.27%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06953'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.18%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.92'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.43%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.029'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.23%  '

$ws.Range("E22").Value = '  +0.16%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.84'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.82%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.392.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.47%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.343'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.96%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.783'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.73%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.75'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.64%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '159.32'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.34%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.713'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -7.25%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '145.32'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.88%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.193'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.86%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.635'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +8.46%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.839.05'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.82%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08227'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.82%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.015'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.59%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.03008'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.46%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.874'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.15%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2776'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.23%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.09470'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.41%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.24'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.70%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.488'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.38%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7831'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.56%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.42'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.44%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.41'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.27%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.558'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.03%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.7041'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.08%  '

$ws.Range("E47").Value = '  -1.43%  '

$ws.Range("E48").Value = '  +3.05%  '

$ws.Range("E49").Value = '  +0.11%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.310'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.70%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '136.80'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.11%  '
